$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.859648466110229
$ws.Range("B1").Value = 4.159233570098877
$ws.Range("C1").Value = 3.501550436019897
$ws.Range("D1").Value = 2.237272977828979
$ws.Range("E1").Value = 1.96209728717804
